$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are stored as text (matches source inlineStr cells),
# not auto-converted to numbers by COM value assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.485.62'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.527.06'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '537.41'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '141.03'
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").Value = '2.538.19'
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").Value = '0.0996'
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").Value = '5.43'
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").Value = '2.971.36'
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("D15").Value = '23.29'
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("D16").Value = '59.335.70'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '2.533.08'
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Value = '11.01'
$ws.Range("E19").Value = '  -2.31%  '
$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").Value = '322.00'
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '5.82'
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("D24").Value = '62.12'
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("D25").Value = '0.421'
$ws.Range("E25").Value = '  -4.36%  '
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  +2.98%  '
$ws.Range("D27").Value = '0.994'
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").Value = '7.86'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("D29").Value = '6.81'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").Value = '0.0₃0775'
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").Value = '162.31'
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.15'
$ws.Range("E33").Value = '  -9.06%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").Value = '1.45'
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").Value = '18.50'
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").Value = '4.23'
$ws.Range("E37").Value = '  -4.72%  '
$ws.Range("D38").Value = '1.60'
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("D39").Value = '36.98'
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '5.41'
$ws.Range("E40").Value = '  -8.30%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.66'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = '290.01'
$ws.Range("E42").Value = '  -5.47%  '
$ws.Range("D43").Value = '0.809'
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("D44").Value = '0.997'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '0.602'
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("D47").Value = '124.93'
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").Value = '0.0930'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").Value = '18.67'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").Value = '0.0512'
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("E51").Value = '  -1.72%  '
